# Restructure the nomination summary sheet so each branch's rows carry a
# qualified "<Branch>, <Metric>" label, add an Air Force section (mirroring
# Army/Navy/Marine Corps/Civilian), and rename/reorder the summary totals at
# the bottom of the sheet (e.g. "Total carryover" -> "Total new nominations" /
# "Total carryover nominations" / "Total recess reappointments ").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label/value layout (column A label, column B value or $null to leave blank)
$rows = @(
    @("Labels", "Values"),
    @("Congress", 98),
    @("Session", 2),
    @("Start Date", 30704),
    @("End Date", 30967),
    @("Army ", $null),
    @("     Army, New nominations", 14031),
    @("     Army, Confirmed ", 14031),
    @("Navy ", $null),
    @("     Navy, New nominations", 8855),
    @("     Navy, Confirmed ", 8855),
    @("Air Force", $null),
    @("     Air Force, New nominations", 11818),
    @("     Air Force, Carryover nominations", 26),
    @("     Air Force, Confirmed ", 11844),
    @("Marine Corps ", $null),
    @("     Marine Corps, New nominations", 2995),
    @("     Marine Corps, Confirmed ", 2995),
    @("Civilian ", $null),
    @("     Civilian, New nominations", 4127),
    @("     Civilian, Confirmed ", 4001),
    @("     Civilian, Unconfirmed ", 107),
    @("     Civilian, Withdrawn ", 2),
    @("     Civilian, Recess reappointments ", 17),
    @("Total new nominations", 41826),
    @("Total carryover nominations", 26),
    @("Total confirmed ", 41726),
    @("Total unconfirmed ", 107),
    @("Total withdrawn", 2),
    @("Total recess reappointments ", 17)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $label = $rows[$i][0]
    $value = $rows[$i][1]

    $ws.Cells.Item($r, 1).Value = $label
    if ($null -eq $value) {
        $ws.Cells.Item($r, 2).ClearContents()
    } else {
        $ws.Cells.Item($r, 2).Value = $value
    }
}

# A few rows moved between the "general/right-aligned" number style and the
# "#,##0" thousands-separated style as the totals were reshuffled (e.g. row 25
# used to be blank and is now a #,##0 total; row 27 used to hold the #,##0
# total and is now a plain count). Re-stamp those four cells' formats from
# same-styled neighbours so the existing style entries are reused rather than
# minting new ones.
$ws.Cells.Item(7, 2).Copy() | Out-Null
$ws.Cells.Item(25, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> #,##0
$ws.Cells.Item(27, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> #,##0

$ws.Cells.Item(22, 2).Copy() | Out-Null
$ws.Cells.Item(26, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> General
$ws.Cells.Item(28, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> General

$excel.CutCopyMode = $false

# Re-apply the values after the format paste (PasteSpecial only touches
# formatting, but keep this deterministic regardless of paste semantics).
$ws.Cells.Item(25, 2).Value = 41826
$ws.Cells.Item(26, 2).Value = 26
$ws.Cells.Item(27, 2).Value = 41726
$ws.Cells.Item(28, 2).Value = 107

# The sheet had 31 rows before; it now has 30, so remove the trailing row.
$ws.Rows.Item(31).Delete()
